$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.33
Write-Host "done"
